# Generate Report for Handback
#
# Refreshes the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Correspond Handoff/Handback DateTime" columns for the second tracked file
# (e0757593-76dd-44f0-b232-3cdbe96e3da1.md) in both the zh-cn and de-de
# localization status sheets, as a new handback round completed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview: row 3 (e0757593-76dd-44f0-b232-3cdbe96e3da1.md) latest generate date.
$overview.Range("G3").Value = "2016-09-06 02:53:54"

# zh-cn: row 3 handoff/handback timestamps for the same file.
$zhcn.Range("H3").Value = "2016-09-06 02:53:50"
$zhcn.Range("K3").Value = "2016-09-06 02:54:16"

# de-de: row 3 handoff/handback timestamps for the same file.
$dede.Range("H3").Value = "2016-09-06 02:53:54"
$dede.Range("K3").Value = "2016-09-06 02:54:23"
